# "Generate Report for Archive"
#
# The localization status report is being refreshed: items that were
# previously marked "Ready for handoff" are now shown as "In Translation",
# and the Status/zh-cn/de-de columns are narrowed to fit the new (shorter)
# text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status text -------------------------------------------------
# "Overview" sheet keeps one status cell per language (columns E = zh-cn,
# F = de-de) for each of the two rows of data.
$overview.Range("E2:F3").Value = "In Translation"

# The per-language detail sheets keep the status in column C.
$zhcn.Range("C2:C3").Value = "In Translation"
$dede.Range("C2:C3").Value = "In Translation"

# --- Narrow the status columns to match the new, shorter text --------------
$overview.Range("E1:F1").ColumnWidth = 12.576851526896165
$zhcn.Range("C1").ColumnWidth = 12.576851526896165
$dede.Range("C1").ColumnWidth = 12.576851526896165
